# Updates cryptos list data (prices, volumes, and two row reorderings)
# to match the Wed Jan  3 11:55:30 UTC 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '44.361.85'
$ws.Range("E2").Value = '  -2.33%  '

# Row 3
$ws.Range("D3").Value = '2.328.67'
$ws.Range("E3").Value = '  -2.52%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '326.58'
$ws.Range("E5").Value = '  +2.37%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.57'
$ws.Range("E6").Value = '  -8.17%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.625'
$ws.Range("E7").Value = '  -1.39%  '

# Row 8
$ws.Range("E8").Value = '  +0.04%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.602'
$ws.Range("E9").Value = '  -3.82%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.83'
$ws.Range("E10").Value = '  -5.94%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0905'
$ws.Range("E11").Value = '  -2.56%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '8.24'
$ws.Range("E12").Value = '  -4.68%  '

# Row 13
$ws.Range("E13").Value = '  -1.59%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.959'
$ws.Range("E14").Value = '  -5.34%  '

# Row 15
$ws.Range("D15").Value = '2.684.39'
$ws.Range("E15").Value = '  -2.37%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.10'
$ws.Range("E16").Value = '  -4.94%  '

# Row 17
$ws.Range("D17").Value = '2.349.17'
$ws.Range("E17").Value = '  -1.52%  '

# Row 18
$ws.Range("D18").Value = '44.426.02'
$ws.Range("E18").Value = '  -2.09%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '15.11'
$ws.Range("E19").Value = '  +11.04%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.17'
$ws.Range("E20").Value = '  -4.44%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0000104'
$ws.Range("E21").Value = '  -3.39%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.56'
$ws.Range("E22").Value = '  -0.02%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.40'
$ws.Range("E23").Value = '  -3.21%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '254.60'
$ws.Range("E24").Value = '  -3.53%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.23'
$ws.Range("E25").Value = '  -5.14%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.20%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.05'
$ws.Range("E27").Value = '  -2.31%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.25'
$ws.Range("E28").Value = '  -6.62%  '

# Row 29
$ws.Range("E29").Value = '  -4.25%  '

# Row 30
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0936'
$ws.Range("E30").Value = '  -4.66%  '

# Row 31
$ws.Range("B31").Value = 'EthereumClassic'
$ws.Range("C31").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '21.72'
$ws.Range("E31").Value = '  -4.49%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '36.34'
$ws.Range("E32").Value = '  -8.09%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '164.57'
$ws.Range("E33").Value = '  -4.64%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.80'
$ws.Range("E34").Value = '  -4.71%  '

# Row 35
$ws.Range("E35").Value = '  -2.27%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.17'
$ws.Range("E36").Value = '  +2.80%  '

# Row 37
$ws.Range("E37").Value = '  -3.67%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.61'
$ws.Range("E38").Value = '  -7.13%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.88'
$ws.Range("E39").Value = '  +5.99%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.87'
$ws.Range("E40").Value = '  -6.70%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0347'
$ws.Range("E41").Value = '  -4.80%  '

# Row 42
$ws.Range("D42").Value = '1.859.12'
$ws.Range("E42").Value = '  +12.71%  '

# Row 43
$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.00'
$ws.Range("E43").Value = '  +0.31%  '

# Row 44
$ws.Range("B44").Value = 'BitcoinSV'
$ws.Range("C44").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '93.96'
$ws.Range("E44").Value = '  -5.82%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '67.57'
$ws.Range("E45").Value = '  -5.69%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.223'
$ws.Range("E46").Value = '  -7.80%  '

# Row 47
$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.62'
$ws.Range("E47").Value = '  -7.36%  '

# Row 48
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.87'
$ws.Range("E48").Value = '  +1.05%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '82.22'
$ws.Range("E49").Value = '  -6.68%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.11'
$ws.Range("E50").Value = '  -3.90%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '109.59'
$ws.Range("E51").Value = '  -5.50%  '
